$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.211.61"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "'2.641.39"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'589.19"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("D6").Value = "'174.50"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.519"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "'2.640.31"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "'0.171"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "'4.95"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").Value = "'3.132.09"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").Value = "'0.0000186"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "'72.089.95"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "'25.87"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "'2.718.84"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").Value = "'12.21"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").Value = "'7.96"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'372.31"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").Value = "'2.06"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'71.05"
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("D26").Value = "'4.27"
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("D27").Value = "'9.64"
$ws.Range("E27").Value = "  -3.99%  "
$ws.Range("D28").Value = "'2.779.81"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'0.0₃0953"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "'8.00"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").Value = "'497.51"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("D33").Value = "'1.28"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("D34").Value = "'1.80"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'161.31"
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("D37").Value = "'19.29"
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("D38").Value = "'0.113"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("D39").Value = "'18.89"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "'1.35"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = "  -6.67%  "
$ws.Range("D43").Value = "'2.56"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "'4.90"
$ws.Range("E44").Value = "  -4.27%  "
$ws.Range("D45").Value = "'0.327"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "'152.80"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'0.548"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").Value = "'3.65"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("D50").Value = "'1.68"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("E51").Value = "  -1.76%  "
